# Update "F" column (想去人数 / "want-to-go" head-count) figures to the
# latest scraped numbers, per "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 3084
$ws.Range("F14").Value = 10
$ws.Range("F18").Value = 1946
$ws.Range("F19").Value = 1773
$ws.Range("F22").Value = 1061
$ws.Range("F24").Value = 1749
$ws.Range("F26").Value = 1182
$ws.Range("F28").Value = 41
$ws.Range("F29").Value = 1104
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 2648
$ws.Range("F34").Value = 2956
$ws.Range("F35").Value = 2144
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 34
$ws.Range("F43").Value = 365
$ws.Range("F47").Value = 670
$ws.Range("F49").Value = 51

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 70
$ws.Range("F23").Value = 59

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 1814
$ws.Range("F8").Value  = 2866
$ws.Range("F12").Value = 396
$ws.Range("F14").Value = 8029

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 1814
$ws.Range("F11").Value = 396
$ws.Range("F15").Value = 70
$ws.Range("F17").Value = 1946
$ws.Range("F19").Value = 1061
$ws.Range("F21").Value = 1749
$ws.Range("F23").Value = 1182
$ws.Range("F26").Value = 1104
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 2648
$ws.Range("F34").Value = 2956
$ws.Range("F35").Value = 2144
$ws.Range("F41").Value = 34
$ws.Range("F44").Value = 59
